$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.590.26"
$ws.Range("E2").Value = "  +1.28%  "

$ws.Range("D3").Value = "2.250.87"
$ws.Range("E3").Value = "  +0.58%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'246.37"
$ws.Range("E5").Value = "  -0.04%  "

$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").Value = "'76.48"
$ws.Range("E7").Value = "  +0.92%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "'0.626"
$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("D10").Value = "'44.55"
$ws.Range("E10").Value = "  +10.44%  "

$ws.Range("D11").Value = "'0.0955"
$ws.Range("E11").Value = "  +0.36%  "

$ws.Range("E12").Value = "  +3.40%  "

$ws.Range("E13").Value = "  -0.32%  "

$ws.Range("E14").Value = "  -0.68%  "

$ws.Range("D15").Value = "'0.865"
$ws.Range("E15").Value = "  +0.89%  "

$ws.Range("D16").Value = "2.269.57"
$ws.Range("E16").Value = "  +0.64%  "

$ws.Range("D17").Value = "42.392.87"
$ws.Range("E17").Value = "  +1.17%  "

$ws.Range("E18").Value = "  +3.75%  "

$ws.Range("D19").Value = "'6.21"
$ws.Range("E19").Value = "  +1.29%  "

$ws.Range("D20").Value = "'72.28"
$ws.Range("E20").Value = "  +0.92%  "

$ws.Range("D21").Value = "'11.01"
$ws.Range("E21").Value = "  +53.14%  "

$ws.Range("E22").Value = "  -0.48%  "

$ws.Range("D23").Value = "'230.97"
$ws.Range("E23").Value = "  -0.17%  "

$ws.Range("D24").Value = "'11.81"
$ws.Range("E24").Value = "  +3.01%  "

$ws.Range("E27").Value = "  +0.89%  "

$ws.Range("E28").Value = "  +4.02%  "

$ws.Range("D29").Value = "'167.45"
$ws.Range("E29").Value = "  -0.95%  "

$ws.Range("D30").Value = "'20.75"
$ws.Range("E30").Value = "  +1.02%  "

$ws.Range("D31").Value = "'5.85"
$ws.Range("E31").Value = "  +20.60%  "

$ws.Range("E32").Value = "  -1.82%  "

$ws.Range("D33").Value = "'31.77"
$ws.Range("E33").Value = "  -4.96%  "

$ws.Range("E34").Value = "  -0.79%  "

$ws.Range("D35").Value = "'4.80"
$ws.Range("E35").Value = "  +6.31%  "

$ws.Range("E36").Value = "  +0.07%  "

$ws.Range("D37").Value = "'0.0318"
$ws.Range("E37").Value = "  +5.93%  "

$ws.Range("D38").Value = "'14.23"
$ws.Range("E38").Value = "  +6.45%  "

$ws.Range("E39").Value = "  +1.05%  "

$ws.Range("D40").Value = "'5.80"
$ws.Range("E40").Value = "  -2.15%  "

$ws.Range("D41").Value = "'64.41"
$ws.Range("E41").Value = "  +6.51%  "

$ws.Range("E42").Value = "  +0.54%  "

$ws.Range("D43").Value = "'108.37"
$ws.Range("E43").Value = "  -3.17%  "

$ws.Range("D44").Value = "'8.92"
$ws.Range("E44").Value = "  +2.32%  "

$ws.Range("E45").Value = "  +2.25%  "

$ws.Range("D46").Value = "'0.996"
$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("D47").Value = "'2.42"
$ws.Range("E47").Value = "  +7.96%  "

$ws.Range("E48").Value = "  +0.70%  "

$ws.Range("D49").Value = "'1.20"
$ws.Range("E49").Value = "  +2.23%  "

$ws.Range("D50").Value = "'4.17"
$ws.Range("E50").Value = "  -0.09%  "

$ws.Range("E51").Value = "  +0.98%  "
